$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in A1 (08:04 -> 09:04) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 09:04"

# --- Swap order of "San Marino" (row 114) and "Georgia" (row 115) ---
# Row 114 becomes Georgia, row 115 becomes San Marino (names swap, data stays attached to the country)
$ws.Cells.Item(114, 1).Value = "Georgia"
$ws.Cells.Item(115, 1).Value = "San Marino"

# --- Swap order of "Belice" (row 192) and "Nueva Caledonia" (row 193) ---
$ws.Cells.Item(192, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(193, 1).Value = "Belice"

# --- Update numeric stats (countries & provincias Spain update) ---

# Row 38: Ucrania
$ws.Cells.Item(38, 2).Value = 14710
$ws.Cells.Item(38, 3).Value = 515
$ws.Cells.Item(38, 4).Value = 2909
$ws.Cells.Item(38, 5).Value = 11425
$ws.Cells.Item(38, 6).Value = 190
$ws.Cells.Item(38, 7).Value = 15
$ws.Cells.Item(38, 8).Value = 376

# Row 49: Chequia
$ws.Cells.Item(49, 2).Value = 8078
$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(49, 4).Value = 4414

# Row 59: Kazajistan
$ws.Cells.Item(59, 2).Value = 4922
$ws.Cells.Item(59, 3).Value = 88
$ws.Cells.Item(59, 5).Value = 3260

# Row 114: now Georgia's data
$ws.Cells.Item(114, 2).Value = 626
$ws.Cells.Item(114, 3).Value = 3
$ws.Cells.Item(114, 4).Value = 297
$ws.Cells.Item(114, 5).Value = 319
$ws.Cells.Item(114, 6).Value = 6
$ws.Cells.Item(114, 8).Value = 10

# Row 115: now San Marino's data
$ws.Cells.Item(115, 4).Value = 114
$ws.Cells.Item(115, 5).Value = 468
$ws.Cells.Item(115, 6).Value = 4
$ws.Cells.Item(115, 8).Value = 41

# Row 125: Taiwan
$ws.Cells.Item(125, 4).Value = 361
$ws.Cells.Item(125, 5).Value = 73

# Row 192: now Nueva Caledonia's data
$ws.Cells.Item(192, 4).Value = 18
$ws.Cells.Item(192, 8).Value = 0

# Row 193: now Belice's data
$ws.Cells.Item(193, 4).Value = 16
$ws.Cells.Item(193, 8).Value = 2
